$p = $ppt.ActivePresentation

# --- Delete the "Web" cloud shape and its connector from the single slide ---
$s = $p.Slides.Item(1)
$s.Shapes.Item("Elbow Connector 51").Delete()
$s.Shapes.Item("Cloud 50").Delete()

# --- Refresh the cached "datetimeFigureOut" field text (slide master,
#     slide layouts, notes master) from 12/15/2018 to 5/29/2019 ---
$newDate = "5/29/2019"

function Update-DateField {
    param($shapes)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tf = $shp.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                if ($tr.Text -eq "12/15/2018") {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

Update-DateField $p.SlideMaster.Shapes
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    Update-DateField $p.SlideMaster.CustomLayouts.Item($li).Shapes
}
Update-DateField $p.NotesMaster.Shapes
